$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Respondent ratio (%)" row entirely (old row 3); the rows below
#    (Revisit count, Fixation count, Dwell time (ms)/(%), Fixation duration (ms),
#    First fixation duration (ms)) shift up by one.
$ws.Rows.Item(3).Delete()

# 2. Correct the data-cleaning recalculation for the shifted metric rows.
#    "First fixation duration (ms)" (now row 8) keeps its original values unchanged.
# Revisit count (now row 3)
$ws.Range("E3").Value = 38
$ws.Range("F3").Value = 39
$ws.Range("H3").Value = 12
$ws.Range("I3").Value = 18
$ws.Range("V3").Value = 30
# Fixation count (now row 4)
$ws.Range("E4").Value = 125
$ws.Range("F4").Value = 104
$ws.Range("H4").Value = 35
$ws.Range("I4").Value = 30
$ws.Range("V4").Value = 242
# Dwell time (ms) (now row 5)
$ws.Range("E5").Value = 36993.5
$ws.Range("F5").Value = 33204.9
$ws.Range("H5").Value = 13084.27
$ws.Range("I5").Value = 10196.96
$ws.Range("V5").Value = 77943.66
# Dwell time (%) (now row 6)
$ws.Range("B6").Value = 0.1
$ws.Range("E6").Value = 25.41
$ws.Range("F6").Value = 22.81
$ws.Range("G6").Value = 0.52
$ws.Range("H6").Value = 8.99
$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 6.59
$ws.Range("K6").Value = 0.5
$ws.Range("M6").Value = 2.7
$ws.Range("N6").Value = 0.91
$ws.Range("S6").Value = 0.38
$ws.Range("T6").Value = 0.22
$ws.Range("V6").Value = 53.54
$ws.Range("Y6").Value = 0.74
# Fixation duration (ms) (now row 7)
$ws.Range("E7").Value = 295.95
$ws.Range("F7").Value = 319.28
$ws.Range("H7").Value = 373.84
$ws.Range("I7").Value = 339.9
$ws.Range("V7").Value = 322.08

# 3. Remove the two trailing blank filler rows (old rows 10-11, now rows 9-10).
$ws.Range("A9:AA10").Delete()

# 4. Strip the bold/centered/bordered header style from row 1, and clear the
#    leftover "Unnamed: 0" label in A1 so the header row is plain again.
$ws.Range("A1:AA1").ClearFormats()
$ws.Range("A1").Value = ""
